$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.47
$ws.Range("I2").Value = 5.8
$ws.Range("L2").Value = 5.4
$ws.Range("N2").Value = 11.8
$ws.Range("P2").Value = 3.88
$ws.Range("Q2").Value = 1.62
$ws.Range("R2").Value = 2.02
$ws.Range("X2").Value = 6.3
$ws.Range("AE2").Value = 12.5
$ws.Range("AF2").Value = 45
$ws.Range("AK2").Value = 90
$ws.Range("AO2").Value = 7
$ws.Range("AS2").Value = 175
$ws.Range("AX2").Value = 32
